$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45202 -> 45203, i.e. 2023-10-03 -> 2023-10-04) for every data row (2-39).
for ($r = 2; $r -le 39; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45202) {
        $cell.Value2 = 45203
    }
}
